# Sanity Test Case - Registration Flow is added now
#
# - Mark the 5 "No" Runmode cells on "Test Cases" as "YES" (these are the
#   rows for MA_AccountEdit1, MA_FulfillmentTool, HomeScreen1, HomeScreen2,
#   HomeScreen3) so the new tests actually run.
# - Make "Test Cases" the active/selected sheet & tab, with D3 selected,
#   instead of SlideTool2 / B21.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("D2").Value = "YES"
$ws.Range("D3").Value = "YES"
$ws.Range("D4").Value = "YES"
$ws.Range("D5").Value = "YES"
$ws.Range("D6").Value = "YES"

$ws.Activate()
$ws.Range("D3").Select()
